# Update the "Avverkningsanmälningar" sheet:
#  1. Every existing data row's "Förändrad" (C) date bumps from 45184 to 45186.
#  2. Every HYPERLINK(...) formula in columns S:Y gains a second (friendly
#     display text) argument equal to the link's filename without extension.
#  3. Row 423 picks up an explicit row height (like every other data row).
#  4. A new case, "A 43498-2023", is appended as row 424.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bump the "Förändrad" column for every existing data row ------------
$ws.Range("C2:C423").Value = 45186

# --- 2. Add the friendly-name argument to every HYPERLINK formula ----------
for ($r = 2; $r -le 423; $r++) {
    for ($c = 19; $c -le 25; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        if ($cell.HasFormula) {
            $f = $cell.Formula
            if ($f -match '^=HYPERLINK\("([^"]+)"\)$') {
                $url = $matches[1]
                $parts = $url.Split("/")
                $last = $parts[$parts.Length - 1]
                $dot = $last.LastIndexOf(".")
                $base = $last.Substring(0, $dot)
                $cell.Formula = '=HYPERLINK("' + $url + '", "' + $base + '")'
            }
        }
    }
}

# --- 3. Row 423 now carries an explicit row height --------------------------
$ws.Rows.Item(423).RowHeight = 15

# --- 4. Append the new case as row 424 --------------------------------------
$dateFormat = $ws.Range("B423").NumberFormat()

$ws.Range("A424").Value = "A 43498-2023"

$ws.Range("B424").Value = 45184
$ws.Range("B424").NumberFormat = $dateFormat

$ws.Range("C424").Value = 45186
$ws.Range("C424").NumberFormat = $dateFormat

$ws.Range("D424").Value = "DALARNAS LÄN"
$ws.Range("E424").Value = "MORA"
$ws.Range("F424").Value = "Bergvik skog väst AB"
$ws.Range("G424").Value = 2.2
$ws.Range("H424").Value = 0
$ws.Range("I424").Value = 0
$ws.Range("J424").Value = 0
$ws.Range("K424").Value = 0
$ws.Range("L424").Value = 0
$ws.Range("M424").Value = 0
$ws.Range("N424").Value = 0
$ws.Range("O424").Value = 0
$ws.Range("P424").Value = 0
$ws.Range("Q424").Value = 0

$ws.Range("R424").Value = ""
$ws.Range("R424").WrapText = $true
